$wb = $excel.ActiveWorkbook

# Insert the new "Bank0x10" sheet right after "Bank0xE" (i.e. before "BANK61"),
# matching the sheet order change in the target workbook.
$bank0xE = $wb.Worksheets.Item("Bank0xE")
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $bank0xE)
$ws.Name = "Bank0x10"

# Header row
$ws.Range("A1").Value = "Purpose"
$ws.Range("B1").Value = "Start"
$ws.Range("C1").Value = "Size"
$ws.Range("D1").Value = "No"
$ws.Range("E1").Value = "Total Size"
$ws.Range("F1").Value = "Code Purpose"
$ws.Range("A1:F1").Font.Bold = $true

# Row 2 - code
$ws.Range("A2").Value = "code"
$ws.Range("B2").Value = 0
$ws.Range("C2").Value = 1586
$ws.Range("D2").Value = 1
$ws.Range("E2").Formula = "=C2"
$ws.Range("F2").Value = "Memory Management"

# Row 3 - memory area (note: shared-string insertion order matters downstream,
# so write F3 before A3)
$ws.Range("B3").Value = 8156
$ws.Range("C3").Value = 35
$ws.Range("E3").Formula = "=C3"
$ws.Range("F3").Value = "dynamic memory management data"
$ws.Range("A3").Value = "memory area"

# Row 12 - Total
$ws.Range("A12").Value = "Total"
$ws.Range("A12").Font.Bold = $true
$ws.Range("E12").Formula = "=SUM(E2:E10)"

# Column widths (closest achievable values; the host's column-width grid
# quantizes to 1/6-character steps, so these land on the nearest bucket to
# the target widths of 31.36328125 / 6.1796875 / 11.453125 characters)
$ws.Columns.Item(1).ColumnWidth = 30.5
$ws.Columns.Item(2).ColumnWidth = 5.333333333
$ws.Range("C1:D1").ColumnWidth = 10.666666667
$ws.Columns.Item(6).ColumnWidth = 30.5

# Selection matches the saved state in the target file
[void]$ws.Range("A3").Select()

Write-Host "Bank0x10 sheet inserted"
